$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add Task 7 / Task 8 / Task 9 rows with their descriptions
$ws.Cells.Item(9, 1).Value2  = "Task 7"
$ws.Cells.Item(9, 3).Value2  = "Create hedge portfolio from factors"

$ws.Cells.Item(10, 1).Value2 = "Task 8 "
$ws.Cells.Item(10, 3).Value2 = "Run CAPM regression on hedge portfolio to check for excess return"

$ws.Cells.Item(11, 1).Value2 = "Task 9"
$ws.Cells.Item(11, 3).Value2 = "Re-Create Table 3 and Table 4 from Hanauer, Lauterbach Paper"

# Update the active selection to match the saved view state
$ws.Range("B8").Select()

$wb.Save()
